# Updates the cryptocurrency price/volume table (columns D and E, rows 2-51)
# on Sheet1 to match the refreshed GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "27.583.81"
$ws.Range("E2").Value2 = "  +2.47%  "
$ws.Range("D3").Value2 = "1.851.69"
$ws.Range("E3").Value2 = "  +2.00%  "
$ws.Range("D4").Value2 = "'1.031"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value2 = "  +2.55%  "
$ws.Range("D5").Value2 = "'321.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value2 = "  +3.40%  "
$ws.Range("D6").Value2 = "'1.030"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value2 = "  +2.44%  "
$ws.Range("D7").Value2 = "'0.4382"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value2 = "  +2.16%  "
$ws.Range("D8").Value2 = "'0.3769"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value2 = "  +2.17%  "
$ws.Range("D9").Value2 = "'0.07403"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value2 = "  +2.28%  "
$ws.Range("D10").Value2 = "'0.8749"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value2 = "  +1.50%  "
$ws.Range("D11").Value2 = "'21.48"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value2 = "  +2.09%  "
$ws.Range("D12").Value2 = "1.863.73"
$ws.Range("E12").Value2 = "  -7.42%  "
$ws.Range("D13").Value2 = "'5.525"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value2 = "  +2.53%  "
$ws.Range("D14").Value2 = "'6.687"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value2 = "  +0.72%  "
$ws.Range("D15").Value2 = "'0.07220"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value2 = "  +4.49%  "
$ws.Range("D16").Value2 = "'82.74"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value2 = "  +2.63%  "
$ws.Range("D17").Value2 = "'1.037"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value2 = "  +3.05%  "
$ws.Range("D18").Value2 = "'0.000009047"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value2 = "  +1.41%  "
$ws.Range("D19").Value2 = "'1.030"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value2 = "  +2.41%  "
$ws.Range("D20").Value2 = "'15.41"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value2 = "  +1.56%  "
$ws.Range("D21").Value2 = "27.599.63"
$ws.Range("E21").Value2 = "  +2.34%  "
$ws.Range("D22").Value2 = "'5.252"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value2 = "  +1.35%  "
$ws.Range("D23").Value2 = "'11.35"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value2 = "  +2.61%  "
$ws.Range("D24").Value2 = "2.075.02"
$ws.Range("E24").Value2 = "  -7.02%  "
$ws.Range("D25").Value2 = "'157.73"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value2 = "  +2.53%  "
$ws.Range("D26").Value2 = "'1.927"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value2 = "  +2.26%  "
$ws.Range("D27").Value2 = "'18.70"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value2 = "  +2.50%  "
$ws.Range("D28").Value2 = "'5.272"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value2 = "  +1.43%  "
$ws.Range("D29").Value2 = "'1.962"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value2 = "  +4.55%  "
$ws.Range("D30").Value2 = "'116.85"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value2 = "  +1.51%  "
$ws.Range("D31").Value2 = "'0.09039"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value2 = "  +0.92%  "
$ws.Range("D32").Value2 = "'0.7625"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value2 = "  +2.69%  "
$ws.Range("D33").Value2 = "'1.194"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value2 = "  +2.80%  "
$ws.Range("D34").Value2 = "'4.500"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value2 = "  +1.82%  "
$ws.Range("D35").Value2 = "'2.886"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value2 = "  +2.98%  "
$ws.Range("D36").Value2 = "'1.031"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value2 = "  +1.91%  "
$ws.Range("E37").Value2 = "  +2.27%  "
$ws.Range("D38").Value2 = "'0.01975"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value2 = "  +2.89%  "
$ws.Range("D39").Value2 = "'0.05293"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value2 = "  +1.35%  "
$ws.Range("D40").Value2 = "'0.5154"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value2 = "  +1.69%  "
$ws.Range("D41").Value2 = "'2.804"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value2 = "  +2.47%  "
$ws.Range("D42").Value2 = "'0.1673"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value2 = "  +1.88%  "
$ws.Range("D43").Value2 = "'6.718"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value2 = "  +4.61%  "
$ws.Range("D44").Value2 = "'8.474"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value2 = "  +2.95%  "
$ws.Range("D45").Value2 = "'108.94"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value2 = "  +1.95%  "
$ws.Range("E46").Value2 = "  +0.90%  "
$ws.Range("D47").Value2 = "'1.709"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value2 = "  +3.23%  "
$ws.Range("D48").Value2 = "'0.06400"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value2 = "  +1.56%  "
$ws.Range("D49").Value2 = "'0.4644"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value2 = "  +1.98%  "
$ws.Range("D50").Value2 = "'1.863"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value2 = "  +3.58%  "
$ws.Range("D51").Value2 = "'39.21"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value2 = "  +4.36%  "
